# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" sheet at the front of the workbook with
#    ID / NAME / BATTING_HAND / BOWL_STYLE columns for player 4313.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, and replace the full
#    howstat.com scorecard URL values with just the bare numeric match
#    code that used to be the `MatchCode=` query parameter.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, inserted before everything else -------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'4313"
$playerInfo.Range("B2").Value = "Sharafuddin Ashraf"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- 2. MATCH_CARD_LINK -> MATCH_CODE on the remaining sheets ----------
foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Player Info") {
        continue
    }

    # Find the MATCH_CARD_LINK header cell in row 1 and note its column.
    $linkCol = 0
    for ($c = 1; $c -le 20; $c++) {
        $headerCell = $ws.Cells.Item(1, $c)
        $headerVal = $headerCell.Value2
        if ($headerVal -eq "MATCH_CARD_LINK") {
            $linkCol = $c
        }
    }

    if ($linkCol -eq 0) {
        continue
    }

    $ws.Cells.Item(1, $linkCol).Value = "MATCH_CODE"

    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $linkCol)
        $url = $cell.Value2
        if ($url -match "MatchCode=(\d+)") {
            $cell.Value = "'" + $matches[1]
        }
    }
}
